{"js": "// Replace the three-digit x one-digit multiplication facts in the\n// answer table with the new values from the updated worksheet.\nconst replacements = [\n  [\"108\u00d78=864\", \"324\u00d74=1296\"],\n  [\"809\u00d75=4045\", \"748\u00d75=3740\"],\n  [\"378\u00d77=2646\", \"679\u00d76=4074\"],\n  [\"530\u00d77=3710\", \"572\u00d74=2288\"],\n  [\"299\u00d77=2093\", \"707\u00d75=3535\"],\n  [\"421\u00d74=1684\", \"239\u00d76=1434\"],\n  [\"225\u00d76=1350\", \"370\u00d73=1110\"],\n  [\"322\u00d75=1610\", \"529\u00d73=1587\"],\n  [\"889\u00d76=5334\", \"376\u00d76=2256\"],\n  [\"837\u00d79=7533\", \"443\u00d74=1772\"],\n  [\"635\u00d76=3810\", \"194\u00d79=1746\"],\n  [\"336\u00d72=672\", \"610\u00d75=3050\"],\n  [\"221\u00d77=1547\", \"633\u00d73=1899\"],\n  [\"958\u00d79=8622\", \"670\u00d74=2680\"],\n  [\"227\u00d78=1816\", \"772\u00d74=3088\"],\n  [\"689\u00d78=5512\", \"394\u00d72=788\"],\n  [\"673\u00d76=4038\", \"431\u00d77=3017\"],\n  [\"522\u00d74=2088\", \"877\u00d72=1754\"],\n  [\"116\u00d78=928\", \"594\u00d72=1188\"],\n  [\"345\u00d72=690\", \"676\u00d76=4056\"],\n  [\"537\u00d77=3759\", \"952\u00d79=8568\"],\n  [\"113\u00d74=452\", \"892\u00d78=7136\"],\n  [\"112\u00d72=224\", \"697\u00d77=4879\"],\n  [\"679\u00d73=2037\", \"495\u00d78=3960\"],\n  [\"849\u00d76=5094\", \"822\u00d74=3288\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit x one-digit multiplication facts in the\n# answer table with the new values from the updated worksheet.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"108\u00d78=864\", \"324\u00d74=1296\"),\n  @(\"809\u00d75=4045\", \"748\u00d75=3740\"),\n  @(\"378\u00d77=2646\", \"679\u00d76=4074\"),\n  @(\"530\u00d77=3710\", \"572\u00d74=2288\"),\n  @(\"299\u00d77=2093\", \"707\u00d75=3535\"),\n  @(\"421\u00d74=1684\", \"239\u00d76=1434\"),\n  @(\"225\u00d76=1350\", \"370\u00d73=1110\"),\n  @(\"322\u00d75=1610\", \"529\u00d73=1587\"),\n  @(\"889\u00d76=5334\", \"376\u00d76=2256\"),\n  @(\"837\u00d79=7533\", \"443\u00d74=1772\"),\n  @(\"635\u00d76=3810\", \"194\u00d79=1746\"),\n  @(\"336\u00d72=672\", \"610\u00d75=3050\"),\n  @(\"221\u00d77=1547\", \"633\u00d73=1899\"),\n  @(\"958\u00d79=8622\", \"670\u00d74=2680\"),\n  @(\"227\u00d78=1816\", \"772\u00d74=3088\"),\n  @(\"689\u00d78=5512\", \"394\u00d72=788\"),\n  @(\"673\u00d76=4038\", \"431\u00d77=3017\"),\n  @(\"522\u00d74=2088\", \"877\u00d72=1754\"),\n  @(\"116\u00d78=928\", \"594\u00d72=1188\"),\n  @(\"345\u00d72=690\", \"676\u00d76=4056\"),\n  @(\"537\u00d77=3759\", \"952\u00d79=8568\"),\n  @(\"113\u00d74=452\", \"892\u00d78=7136\"),\n  @(\"112\u00d72=224\", \"697\u00d77=4879\"),\n  @(\"679\u00d73=2037\", \"495\u00d78=3960\"),\n  @(\"849\u00d76=5094\", \"822\u00d74=3288\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
